$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "‘파이썬 라이브러리를 활용한 머신러닝(개정 2판)’ 사이킷런 1.2.1에서 재실행 완료"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/03/05/%ed%8c%8c%ec%9d%b4%ec%8d%ac-%eb%9d%bc%ec%9d%b4%eb%b8%8c%eb%9f%ac%eb%a6%ac%eb%a5%bc-%ed%99%9c%ec%9a%a9%ed%95%9c-%eb%a8%b8%ec%8b%a0%eb%9f%ac%eb%8b%9d%ea%b0%9c%ec%a0%95-2%ed%8c%90-%ec%82%ac/"

$ws.Range("D28").Value = "Behavior Tree :: Subsumption Architecture, Teleo-Reactive, Decision Trees"
$ws.Range("E28").Value = "https://ropiens.tistory.com/212"

$ws.Range("D32").Value = "Simpson's paradox (심슨의 역설)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/407"
